$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Relocate the document's single "_GoBack" bookmark (Word only ever keeps
#    one) from its current spot - between the "npm," run and the
#    " Node.js, Express, JSON" run - to right after the " Hawkins " run in
#    the very first paragraph (the name line), i.e. immediately before that
#    paragraph's mark. Re-adding a bookmark under the same name moves it, so
#    this single operation both removes it from its old location (goal 3)
#    and places it at the new one (goal 1).
#
#    The engine mis-places a bookmark that is collapsed exactly at a
#    paragraph's content-end boundary, so we temporarily insert a placeholder
#    character after that boundary, add the bookmark in what is now a
#    *mid-paragraph* position (which the engine handles correctly, leaving
#    existing runs untouched), and then remove the placeholder.
# ---------------------------------------------------------------------------
$nameRange = $d.Paragraphs.Item(1).Range
$nameRange.MoveEnd(1, -1)
$boundary = $nameRange.End

$placeholder = $d.Range($boundary, $boundary)
$placeholder.InsertAfter("X")

$bmSpot = $d.Range($boundary, $boundary)
$d.Bookmarks.Add("_GoBack", $bmSpot)

$d.Range($boundary, $boundary + 1).Delete()

# ---------------------------------------------------------------------------
# 2) Bump the spacing-after of the "Find out more about me at..." paragraph
#    from 0 to 120 twips (6 pt).
# ---------------------------------------------------------------------------
$d.Paragraphs.Item(4).SpaceAfter = 6

# ---------------------------------------------------------------------------
# 4) Clear out the final "Microsoft Office: Excel, Word, PowerPoint"
#    paragraph, leaving it empty (the paragraph itself stays, but loses all
#    of its runs/text).
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$lastPara.MoveEnd(1, -1)
$lastPara.Text = ""

Write-Output "done"
